# NIT-9010630901.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The account-statement table previously listed 3 overdue periods (2507,
# 2506, 2505) for the same worker, with "VALOR MORA" / "Cant. Periodos"
# summarizing all three. This edit replaces that with just the first new
# period (2507): the other two period rows are removed and the summary
# fields are corrected to match a single period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# VALOR MORA (E11): was the sum across 3 periods (170820 = 3 * 56940),
# now reflects just the 1 remaining period.
$ws.Range("E11").Value = 56940

# Cant. Periodos (F13): was 3 periods, now 1.
$ws.Range("F13").Value = 1

# Remove the two extra period rows (2506 in row 17, 2505 in row 18) from
# the detail table, leaving only the 2507 row (16). Remaining rows below
# shift up to fill the gap.
$ws.Rows("17:18").Delete()
